$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 71
$ws.Range("D2").Value = 0.01333
$ws.Range("E2").Value = 0.00667
$ws.Range("F2").Value = 1.34113

# Row 3
$ws.Range("C3").Value = 49.25
$ws.Range("D3").Value = 0.04155
$ws.Range("E3").Value = 0.02078
$ws.Range("F3").Value = 2.65682

# Row 4
$ws.Range("C4").Value = 27.95
$ws.Range("D4").Value = 0.06429
$ws.Range("E4").Value = 0.01607
$ws.Range("F4").Value = 2.5546

# Row 5
$ws.Range("C5").Value = 37.7
$ws.Range("D5").Value = 0.0349
$ws.Range("E5").Value = 0.00873
$ws.Range("F5").Value = 2.03793

# Row 6
$ws.Range("C6").Value = 22.8
$ws.Range("D6").Value = 0.03905
$ws.Range("E6").Value = 0.00651
$ws.Range("F6").Value = 1.79519

# Row 7
$ws.Range("C7").Value = 21.95
$ws.Range("D7").Value = 0.08359999999999999
$ws.Range("E7").Value = 0.01393
$ws.Range("F7").Value = 2.8793

# Row 8
$ws.Range("C8").Value = 14.65
$ws.Range("D8").Value = 0.10446
$ws.Range("E8").Value = 0.01306
$ws.Range("F8").Value = 2.69742

# Row 9
$ws.Range("C9").Value = 17.95
$ws.Range("D9").Value = 0.0536
$ws.Range("E9").Value = 0.0067
$ws.Range("F9").Value = 1.95145

# Row 10
$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 0.04352
$ws.Range("E10").Value = 0.00435
$ws.Range("F10").Value = 1.40021

# Row 11
$ws.Range("C11").Value = 9.550000000000001
$ws.Range("D11").Value = 0.1274
$ws.Range("E11").Value = 0.01274
$ws.Range("F11").Value = 2.40263
